# Weekly refresh of the "Cilantro" price series: a new week of data is
# inserted at the top of the date-ordered block (rows 78-79), every
# following week shifts down by one pair of rows, and the oldest week
# (which used to be the last pair, rows 124-125) is appended as a new
# pair of rows at the bottom (126-127).
#
# Only the "Fecha" (D) and "Origen" (O) columns actually differ from
# week to week for this market/product combination - everything else in
# a pair of rows (Primera/Segunda quality rows) is constant - so shifting
# whole rows down reproduces the data faithfully.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) The oldest existing week (rows 124-125) falls off the end of the
#    shifted block and becomes the new last week (rows 126-127).
$ws.Range("A124:R125").Copy($ws.Range("A126:R127"))

# 2) Shift every week from rows 78-123 down by one pair of rows, to make
#    room for the newly-arrived week at rows 78-79.
$ws.Range("A78:R123").Copy($ws.Range("A80:R125"))

# 3) Write the newly-arrived week's date into the now-vacated rows 78-79
#    (the rest of the row - region, quality, prices, etc. - is unchanged
#    from what was already there).
$ws.Range("D78").Value = 44460
$ws.Range("D79").Value = 44460
